$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.497
$ws.Range("D4").Value = -7.787000000000001
$ws.Range("A6").Value = -22.303
$ws.Range("A7").Value = -20.115
$ws.Range("A8").Value = -22.18
$ws.Range("D8").Value = -8.700999999999999
$ws.Range("D9").Value = -7.678
$ws.Range("D12").Value = -6.697
$ws.Range("A16").Value = -22.107
$ws.Range("D17").Value = -8.509
$ws.Range("D18").Value = -8.620000000000001
$ws.Range("D19").Value = -8.047999999999998
$ws.Range("A20").Value = -20.341
$ws.Range("D20").Value = -7.819999999999999
$ws.Range("A21").Value = -19.891
$ws.Range("D26").Value = -7.628
$ws.Range("A28").Value = -21.927
$ws.Range("A29").Value = -21.439
$ws.Range("A30").Value = -21.568
$ws.Range("D31").Value = -7.793000000000001
$ws.Range("A32").Value = -21.713
$ws.Range("D39").Value = -7.699
$ws.Range("A40").Value = -19.937
$ws.Range("D40").Value = -8.16
$ws.Range("D41").Value = -7.928
$ws.Range("D42").Value = -7.958
$ws.Range("D43").Value = -7.878
$ws.Range("A46").Value = -21.801
$ws.Range("D47").Value = -7.569
$ws.Range("D48").Value = -7.717999999999999
$ws.Range("A51").Value = -21.95
$ws.Range("A52").Value = -22.036
$ws.Range("D54").Value = -7.921000000000001
$ws.Range("A57").Value = -22.339
$ws.Range("A59").Value = -22.699
$ws.Range("A62").Value = -22.195
$ws.Range("D62").Value = -8.418000000000001
$ws.Range("D63").Value = -6.876
$ws.Range("D64").Value = -7.068
$ws.Range("A66").Value = -21.615
$ws.Range("A73").Value = -20.082
$ws.Range("A74").Value = -21.256
$ws.Range("D76").Value = -7.891000000000001
$ws.Range("A77").Value = -20.367
$ws.Range("D81").Value = -7.543000000000001
$ws.Range("D84").Value = -8.33
$ws.Range("D89").Value = -8.141999999999999
$ws.Range("A92").Value = -21.566
$ws.Range("D94").Value = -7.267
$ws.Range("A100").Value = -22.383
